$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("sigma_010")
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 28.16930748809891
$ws.Cells.Item(2, 3).Value = 33.04097173550634
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 28.17297950247055
$ws.Cells.Item(3, 3).Value = 33.05692585751255
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 28.17675549325753
$ws.Cells.Item(4, 3).Value = 33.05172290075931
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 28.20040540071667
$ws.Cells.Item(5, 3).Value = 33.03941920188159
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 28.18786957035335
$ws.Cells.Item(6, 3).Value = 33.04891119959652
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 28.16266070185871
$ws.Cells.Item(7, 3).Value = 33.02975921348241
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 28.17712855766201
$ws.Cells.Item(8, 3).Value = 33.08733188744841
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 28.22725258734179
$ws.Cells.Item(9, 3).Value = 33.07582437837945
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 28.17314980835272
$ws.Cells.Item(10, 3).Value = 33.05004013759471
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 28.14708379545179
$ws.Cells.Item(11, 3).Value = 33.03082551432844
$ws.Cells.Item(12, 2).Value = 28.1794592905564
$ws.Cells.Item(12, 3).Value = 33.05117320264898

$ws = $wb.Worksheets.Item("sigma_025")
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 19.85830119490692
$ws.Cells.Item(2, 3).Value = 29.50381495727706
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 19.83298649664522
$ws.Cells.Item(3, 3).Value = 29.46716303048007
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 19.85508146645041
$ws.Cells.Item(4, 3).Value = 29.49877160464699
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 19.84163348472186
$ws.Cells.Item(5, 3).Value = 29.49483906435426
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 19.85184617895172
$ws.Cells.Item(6, 3).Value = 29.45613160740098
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 19.86403301202607
$ws.Cells.Item(7, 3).Value = 29.51284830809138
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 19.83322692685066
$ws.Cells.Item(8, 3).Value = 29.51395367399795
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 19.85860513323254
$ws.Cells.Item(9, 3).Value = 29.48709728076837
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 19.86051346719876
$ws.Cells.Item(10, 3).Value = 29.44687135889752
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 19.82635023505657
$ws.Cells.Item(11, 3).Value = 29.49646257246044
$ws.Cells.Item(12, 2).Value = 19.84825775960407
$ws.Cells.Item(12, 3).Value = 29.4877953458375

$ws = $wb.Worksheets.Item("sigma_050")
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 14.89663588617208
$ws.Cells.Item(2, 3).Value = 25.00601660287861
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 14.88026691438229
$ws.Cells.Item(3, 3).Value = 24.99099851873351
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 14.89171075275804
$ws.Cells.Item(4, 3).Value = 24.96189070191538
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 14.87840419464074
$ws.Cells.Item(5, 3).Value = 24.96007928716441
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 14.89966891085526
$ws.Cells.Item(6, 3).Value = 24.9497176206858
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 14.89375184249389
$ws.Cells.Item(7, 3).Value = 24.98987031777179
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 14.8825148057283
$ws.Cells.Item(8, 3).Value = 24.93601328084555
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 14.86926121449483
$ws.Cells.Item(9, 3).Value = 24.81095106063454
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 14.88608997419887
$ws.Cells.Item(10, 3).Value = 24.92656837590088
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 14.87354539834943
$ws.Cells.Item(11, 3).Value = 24.93933485306041
$ws.Cells.Item(12, 2).Value = 14.88518498940737
$ws.Cells.Item(12, 3).Value = 24.94714406195909
